$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Shift existing data columns A:E to B:F (rows 1-9), from right to left to avoid overwrite
$ws.Range("F1:F9").Value2 = $ws.Range("E1:E9").Value2
$ws.Range("E1:E9").Value2 = $ws.Range("D1:D9").Value2
$ws.Range("D1:D9").Value2 = $ws.Range("C1:C9").Value2
$ws.Range("C1:C9").Value2 = $ws.Range("B1:B9").Value2
$ws.Range("B1:B9").Value2 = $ws.Range("A1:A9").Value2

# New header row
$ws.Range("A1").Value2 = "Metodo"
$ws.Range("B1").Value2 = "Rx"
$ws.Range("C1").Value2 = "Ry"
$ws.Range("D1").Value2 = "CL"
$ws.Range("E1").Value2 = "Entropia"
$ws.Range("F1").Value2 = "SSIM"

# New column A method names
$ws.Range("A2").Value2 = "SMARTER"
$ws.Range("A3").Value2 = "Fuzzy"
$ws.Range("A4").Value2 = "TOPSIS"
$ws.Range("A5").Value2 = "GRA"
$ws.Range("A6").Value2 = "CODAS"
$ws.Range("A7").Value2 = "MABAC"
$ws.Range("A8").Value2 = "VIKOR"
$ws.Range("A9").Value2 = "PROMETHEE II"

# Column widths to match target (columns D/E already width 12/bestFit from source data,
# leave them untouched so that attribute survives; only touch columns that actually change)
$ws.Columns.Item(1).ColumnWidth = 12.5
$ws.Columns.Item(2).ColumnWidth = 3.1666666666666665
$ws.Columns.Item(3).ColumnWidth = 2.333333333333333
$ws.Columns.Item(6).ColumnWidth = 11.166666666666666
